$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, taken from the latest cryptos-list scrape.
# Values that look like plain decimal numbers (e.g. "523.60") get a leading
# apostrophe so Excel stores them as text (matching the source t="inlineStr"
# cells) instead of silently reinterpreting them as numbers and dropping
# trailing zeros (e.g. "523.60" -> 523.6, "61.00" -> 61).
$updates = @(
    @('D2', '60.919.04'),
    @('E2', '  +0.91%  '),
    @('D3', '2.596.73'),
    @('E3', '  +0.47%  '),
    @('E4', '  -0.09%  '),
    @('D5', '523.60'),
    @('E5', '  +3.42%  '),
    @('D6', '154.28'),
    @('E6', '  +0.22%  '),
    @('D7', '0.997'),
    @('E7', '  -0.13%  '),
    @('D8', '0.596'),
    @('E8', '  +2.53%  '),
    @('D9', '6.71'),
    @('E9', '  +2.11%  '),
    @('D10', '0.106'),
    @('E10', '  +2.07%  '),
    @('E11', '  +0.53%  '),
    @('E12', '  +1.57%  '),
    @('D13', '3.048.79'),
    @('E13', '  +0.40%  '),
    @('D14', '60.929.46'),
    @('E14', '  +0.91%  '),
    @('D15', '21.69'),
    @('E15', '  +1.02%  '),
    @('E16', '  +1.04%  '),
    @('D17', '2.598.65'),
    @('E17', '  +0.46%  '),
    @('D18', '4.75'),
    @('E18', '  -1.05%  '),
    @('D19', '353.53'),
    @('E19', '  +2.34%  '),
    @('D20', '10.59'),
    @('E20', '  +1.72%  '),
    @('D21', '6.23'),
    @('E21', '  +1.49%  '),
    @('E22', '  -0.03%  '),
    @('D23', '61.00'),
    @('E23', '  +1.53%  '),
    @('D24', '0.427'),
    @('E24', '  +1.43%  '),
    @('B25', 'Kaspa'),
    @('C25', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'),
    @('D25', '0.166'),
    @('E25', '  -0.24%  '),
    @('B26', 'WrappedeETH'),
    @('C26', 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'),
    @('D26', '2.709.13'),
    @('E26', '  +0.47%  '),
    @('D27', '0.997'),
    @('E27', '  +0.05%  '),
    @('D28', '0.0₃0843'),
    @('E28', '  -0.18%  '),
    @('D29', '7.36'),
    @('E29', '  +0.44%  '),
    @('E30', '  -0.03%  '),
    @('D31', '6.37'),
    @('E31', '  +11.86%  '),
    @('D32', '19.36'),
    @('E32', '  +0.14%  '),
    @('D33', '1.61'),
    @('E33', '  +3.51%  '),
    @('D34', '148.73'),
    @('E34', '  -3.15%  '),
    @('E35', '  +5.14%  '),
    @('D36', '0.945'),
    @('E36', '  +10.26%  '),
    @('E37', '  +1.32%  '),
    @('D38', '1.50'),
    @('E38', '  +1.97%  '),
    @('D39', '3.81'),
    @('E39', '  +1.20%  '),
    @('B40', 'Fetch.AI'),
    @('C40', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'),
    @('D40', '0.849'),
    @('E40', '  +0.31%  '),
    @('B41', 'OKB'),
    @('C41', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'),
    @('D41', '36.48'),
    @('E41', '  +1.65%  '),
    @('D42', '287.16'),
    @('E42', '  -2.96%  '),
    @('E43', '  +1.48%  '),
    @('E44', '  +1.64%  '),
    @('D45', '0.0559'),
    @('E45', '  +0.41%  '),
    @('D46', '0.997'),
    @('E46', '  +0.01%  '),
    @('D47', '19.56'),
    @('E47', '  -1.25%  '),
    @('B48', 'RenderToken'),
    @('C48', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'),
    @('D48', '4.88'),
    @('E48', '  +0.65%  '),
    @('B49', 'VeChain'),
    @('C49', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @('D49', '0.0237'),
    @('E49', '  +1.77%  '),
    @('E50', '  +0.07%  '),
    @('D51', '18.98'),
    @('E51', '  +7.50%  ')
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $looksNumeric = $val -match '^[+-]?(\d+\.)?\d+$'
    if ($looksNumeric) {
        $val = "'" + $val
    }
    $ws.Range($addr).Value = $val
}
